# The document previously placed the cursor (saved as the "_GoBack"
# bookmark) in the final, otherwise-empty paragraph at the end of the
# document. The edit re-types part of the phrase "Korisnik može koristiti
# uređaj" (splitting it into "Korisnik može koris" + "titi uređaj"),
# which leaves Word's "_GoBack" bookmark at that new cursor position
# instead of at the end of the document.

$d = $word.ActiveDocument

# Locate the target sentence.
$rng = $d.Content
$found = $rng.Find.Execute("Korisnik može koristiti uređaj", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Split point is right after "Korisnik može koris" (19 characters in).
$splitPos = $rng.Start + "Korisnik može koris".Length

# Re-type the tail of the sentence so the text is produced by two runs,
# exactly as a live edit at that cursor position would do.
$tailRange = $d.Range($splitPos, $rng.End)
$tailRange.Text = ""
$insPoint = $d.Range($splitPos, $splitPos)
$insPoint.InsertAfter("titi uređaj")

# Move the "_GoBack" bookmark to the new cursor position (this both
# removes it from the trailing empty paragraph and adds it here, since
# bookmark names are unique).
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
